# "1st changes of mifos to finflux"
#
# The "Repayment Schedule" sheet gets a new (empty) column inserted at
# column N, pushing the existing "Late" / blank / "Outstanding" columns
# one slot to the right (N,O,P -> O,P,Q). That sheet also becomes the
# active/selected sheet (with a new active-cell selection), while the
# previously-active "NewLoanInput1" sheet loses its tab selection.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the existing "Late" column (N) -
# this shifts N -> O, O -> P, P -> Q and widens the used range to A1:Q14.
$wsRepay.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and move the selection
# to the new right-hand edge of the table (old F16 -> S8 after edits).
$wsRepay.Select()
$wsRepay.Range("S8").Select()
